$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.248686548991202
$ws.Range("C2").Value = 0.6347598823854241
$ws.Range("E2").Value = 0.4183381963713657
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.002409453222045104
$ws.Range("I2").Value = 1.309755037841484

$ws.Range("B3").Value = 1.1112033983988
$ws.Range("C3").Value = 0.5557285926852842
$ws.Range("E3").Value = 0.3643157765080929
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.002416663441703992
$ws.Range("I3").Value = 1.227970449418066

$ws.Range("B4").Value = 1.027504074215244
$ws.Range("C4").Value = 0.5074601859276981
$ws.Range("E4").Value = 0.3313260668759881
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.00242130394989196
$ws.Range("I4").Value = 1.178427338663525

$ws.Range("B5").Value = 0.9935701704329745
$ws.Range("C5").Value = 0.4878504732234887
$ws.Range("E5").Value = 0.3179235800748899
$ws.Range("F5").Value = 0.3390132514326325
$ws.Range("G5").Value = 0.002423248915043381
$ws.Range("I5").Value = 1.158399573988163

$ws.Range("B6").Value = 0.9879458305915136
$ws.Range("C6").Value = 0.4845977897548437
$ws.Range("E6").Value = 0.3157004664510907
$ws.Range("F6").Value = 0.3366681778241372
$ws.Range("G6").Value = 0.002423575139066763
$ws.Range("I6").Value = 1.155083504747168

$ws.Range("B7").Value = 1.027045731686542
$ws.Range("C7").Value = 0.5071954853426064
$ws.Range("E7").Value = 0.3311451554405238
$ws.Range("F7").Value = 0.3529483938368969
$ws.Range("G7").Value = 0.002421329961689977
$ws.Range("I7").Value = 1.178156593710739

$ws.Range("B8").Value = 1.201130215155899
$ws.Range("C8").Value = 0.6074537094848438
$ws.Range("E8").Value = 0.3996712566773368
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("G8").Value = 0.002411895170084785
$ws.Range("I8").Value = 1.281411322529763

$ws.Range("B9").Value = 1.548474245408499
$ws.Range("C9").Value = 0.8063201359917684
$ws.Range("E9").Value = 0.5356817657126527
$ws.Range("F9").Value = 0.5661985755042025
$ws.Range("G9").Value = 0.002395074961288197
$ws.Range("I9").Value = 1.489594526421399

$ws.Range("B10").Value = 1.807734566465058
$ws.Range("C10").Value = 0.9541288060121929
$ws.Range("E10").Value = 0.6368945614587034
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002383725395556979
$ws.Range("I10").Value = 1.646539119170853

$ws.Range("B11").Value = 1.926660233455721
$ws.Range("C11").Value = 1.021813068870529
$ws.Range("E11").Value = 0.6832836257581931
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002378777472905281
$ws.Range("I11").Value = 1.718918392158542

$ws.Range("B12").Value = 1.971844008173548
$ws.Range("C12").Value = 1.047513223460669
$ws.Range("E12").Value = 0.7009051812976708
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002376934461263149
$ws.Range("I12").Value = 1.746477313822993

$ws.Range("B13").Value = 1.962106118896941
$ws.Range("C13").Value = 1.041975050144401
$ws.Range("E13").Value = 0.6971075333441377
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002377330027502471
$ws.Range("I13").Value = 1.740535171584725

$ws.Range("B14").Value = 1.930374496758645
$ws.Range("C14").Value = 1.023926008953083
$ws.Range("E14").Value = 0.6847322289719244
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002378625234313928
$ws.Range("I14").Value = 1.721182614159972

$ws.Range("B15").Value = 1.910957613996857
$ws.Range("C15").Value = 1.012879680548394
$ws.Range("E15").Value = 0.6771593186840761
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002379422570828986
$ws.Range("I15").Value = 1.709348475597579

$ws.Range("B16").Value = 1.799982970308179
$ws.Range("C16").Value = 0.9497149130620528
$ws.Range("E16").Value = 0.6338703432317203
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002384053058917202
$ws.Range("I16").Value = 1.641829529785724

$ws.Range("B17").Value = 1.732161487516919
$ws.Range("C17").Value = 0.911083342851839
$ws.Range("E17").Value = 0.6074064128751928
$ws.Range("F17").Value = 0.6400460337215605
$ws.Range("G17").Value = 0.002386948605244887
$ws.Range("I17").Value = 1.60066725813212

$ws.Range("B18").Value = 1.693244798965736
$ws.Range("C18").Value = 0.8889050634140858
$ws.Range("E18").Value = 0.5922173912552182
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002388634304581073
$ws.Range("I18").Value = 1.577083921689848

$ws.Range("B19").Value = 1.680083934461265
$ws.Range("C19").Value = 0.8814028534939098
$ws.Range("E19").Value = 0.587080048038672
$ws.Range("F19").Value = 0.6191636801734006
$ws.Range("G19").Value = 0.002389208540870633
$ws.Range("I19").Value = 1.569114560916375

$ws.Range("B20").Value = 1.739371581471175
$ws.Range("C20").Value = 0.9151913965562244
$ws.Range("E20").Value = 0.6102201604676623
$ws.Range("F20").Value = 0.6429339538360921
$ws.Range("G20").Value = 0.002386638274667604
$ws.Range("I20").Value = 1.605039454797122

$ws.Range("B21").Value = 1.939690734511998
$ws.Range("C21").Value = 1.02922551612653
$ws.Range("E21").Value = 0.6883656207839692
$ws.Range("F21").Value = 0.7228739723492197
$ws.Range("G21").Value = 0.002378243970664955
$ws.Range("I21").Value = 1.726862774119013

$ws.Range("B22").Value = 2.071484140054736
$ws.Range("C22").Value = 1.104161734862828
$ws.Range("E22").Value = 0.7397614112366426
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002372936389301032
$ws.Range("I22").Value = 1.807362611314318

$ws.Range("B23").Value = 2.001061136392366
$ws.Range("C23").Value = 1.064127627997095
$ws.Range("E23").Value = 0.7122991975049047
$ws.Range("F23").Value = 0.7472568307916134
$ws.Range("G23").Value = 0.002375752893364089
$ws.Range("I23").Value = 1.764314713290645

$ws.Range("B24").Value = 1.73611166560346
$ws.Range("C24").Value = 0.9133340477011984
$ws.Range("E24").Value = 0.6089479863148455
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002386778509890232
$ws.Range("I24").Value = 1.603062531808291

$ws.Range("B25").Value = 1.453826259618609
$ws.Range("C25").Value = 0.7522465533280638
$ws.Range("E25").Value = 0.4986830249168008
$ws.Range("F25").Value = 0.5279251897347308
$ws.Range("G25").Value = 0.002399446997422821
$ws.Range("I25").Value = 1.432612166527406
